$d = $word.ActiveDocument
$rng = $d.Content
$found = $rng.Find.Execute("go to homepage of ", $true, $false, $false, $false, $false, $true, 1, $false, "go to homepage of ", 2)
Write-Output "found=$found"
